# [Hall Effect] Add figs calibration, update jupyter notebook and fix combined datasets
#
# The Combined_dataset sheet stored R_mn_op (D) and R_no_pm (E) as negative
# resistances; this fixes the sign so the values are reported as their
# magnitudes (positive), matching the corrected data-processing notebook.
# A handful of rows (27-45) also pick up tiny floating point re-computation
# differences (their Hall_Voltage in column C, plus D/E) because the
# upstream notebook recomputed the whole column instead of just negating it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D2" = 1520
    "E2" = 2230
    "D3" = 1290
    "E3" = 1940
    "D4" = 1056
    "E4" = 1560
    "D5" = 895.1
    "E5" = 1350
    "D6" = 777.5
    "E6" = 1171
    "D7" = 688.7
    "E7" = 1034.5
    "D8" = 618.9
    "E8" = 927.5
    "D9" = 563.3
    "E9" = 841.5
    "D10" = 518.3
    "E10" = 771.3
    "D11" = 481.1999999999999
    "E11" = 713.3
    "D12" = 450.9
    "E12" = 665.6999999999999
    "D13" = 426
    "E13" = 626.5999999999999
    "D14" = 404.7
    "E14" = 592.8
    "D15" = 387.2
    "E15" = 564.9
    "D16" = 372.1
    "E16" = 540.6999999999999
    "D17" = 359.92
    "E17" = 521.0799999999999
    "D18" = 349.53
    "E18" = 503.95
    "D19" = 340.85
    "E19" = 489.35
    "D20" = 333.54
    "E20" = 477.19
    "D21" = 327.2699999999999
    "E21" = 466.3
    "D22" = 322.67
    "E22" = 458.04
    "D23" = 318.76
    "E23" = 450.73
    "D24" = 315.4299999999999
    "E24" = 444.26
    "D25" = 313.12
    "E25" = 439.27
    "D26" = 311.4999999999999
    "E26" = 435.1899999999999
    "C27" = 1.6
    "D27" = 310
    "E27" = 431
    "D28" = 308.9999999999999
    "E28" = 428.9999999999999
    "D29" = 310
    "E29" = 427
    "D30" = 311
    "E30" = 426
    "C31" = 1.379999999999999
    "D31" = 312
    "E31" = 426
    "C32" = 1.4
    "D32" = 313.9999999999999
    "E32" = 426
    "C33" = 1.4
    "D33" = 316
    "E33" = 426
    "C34" = 1.300000000000001
    "D34" = 318
    "E34" = 427
    "C35" = 1.299999999999999
    "D35" = 321
    "E35" = 427.9999999999999
    "C36" = 1.300000000000001
    "D36" = 325
    "E36" = 430
    "C37" = 1.300000000000001
    "D37" = 327
    "E37" = 431
    "C38" = 1.34
    "D38" = 330
    "E38" = 433.9999999999999
    "C39" = 1.199999999999999
    "D39" = 333.9999999999999
    "E39" = 435
    "C40" = 1.200000000000001
    "D40" = 337
    "E40" = 437
    "C41" = 1.199999999999999
    "D41" = 342
    "E41" = 438.9999999999999
    "C42" = 1.100000000000001
    "D42" = 347.9999999999999
    "E42" = 442
    "C43" = 1.199999999999999
    "D43" = 352
    "E43" = 443.9999999999999
    "C44" = 1.100000000000001
    "D44" = 356
    "E44" = 445
    "C45" = 1.2
    "D45" = 357.9999999999999
    "E45" = 447
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

Write-Output ("Updated {0} cells" -f $values.Count)
